$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 0.9999970461623668
$ws.Range("E2").Value = 1.003512603051346
$ws.Range("F2").Value = 0.998053937717466
$ws.Range("J2").Value = 1.005312350954268
$ws.Range("L2").Value = 1.006403886417998
$ws.Range("M2").Value = 1.000962071010938
$ws.Range("N2").Value = 1.005586813710726
$ws.Range("C3").Value = 1.002190696077172
$ws.Range("E3").Value = 1.005418250165749
$ws.Range("F3").Value = 1.000925108801209
$ws.Range("J3").Value = 1.007127342732228
$ws.Range("L3").Value = 1.008108944823935
$ws.Range("M3").Value = 1.003628619941815
$ws.Range("N3").Value = 1.006228286426888
$ws.Range("C4").Value = 1.003602290855778
$ws.Range("E4").Value = 1.006644382404553
$ws.Range("F4").Value = 1.002773591153556
$ws.Range("J4").Value = 1.008294031562883
$ws.Range("L4").Value = 1.009204923891011
$ws.Range("M4").Value = 1.005344591826161
$ws.Range("N4").Value = 1.006639578476929
$ws.Range("C5").Value = 1.004193894425285
$ws.Range("E5").Value = 1.007158223888107
$ws.Range("F5").Value = 1.003548520839981
$ws.Range("J5").Value = 1.008782696506916
$ws.Range("L5").Value = 1.009663960881819
$ws.Range("M5").Value = 1.006063784403245
$ws.Range("N5").Value = 1.006811593515251
$ws.Range("C6").Value = 1.004293121253614
$ws.Range("E6").Value = 1.007244405997113
$ws.Range("F6").Value = 1.003678509513616
$ws.Range("J6").Value = 1.00886464045466
$ws.Range("L6").Value = 1.009740935846063
$ws.Range("M6").Value = 1.00618441270587
$ws.Range("N6").Value = 1.006840423675392
$ws.Range("C7").Value = 1.003610203023347
$ws.Range("E7").Value = 1.006651254710765
$ws.Range("F7").Value = 1.002783954247484
$ws.Range("J7").Value = 1.008300568190577
$ws.Range("L7").Value = 1.009211064246
$ws.Range("M7").Value = 1.005354210280487
$ws.Range("N7").Value = 1.006641880437511
$ws.Range("C8").Value = 1.000740052902659
$ws.Range("E8").Value = 1.004158087951371
$ws.Range("F8").Value = 0.9990262421512631
$ws.Range("J8").Value = 1.00592736126651
$ws.Range("L8").Value = 1.00698165440412
$ws.Range("M8").Value = 1.001865241835927
$ws.Range("N8").Value = 1.005804394172819
$ws.Range("C9").Value = 0.9956202764081894
$ws.Range("E9").Value = 0.999709809660648
$ws.Range("F9").Value = 0.9923299142558757
$ws.Range("J9").Value = 1.001684422044947
$ws.Range("L9").Value = 1.002995500289087
$ws.Range("M9").Value = 0.9956418554368391
$ws.Range("N9").Value = 1.004299033365596
$ws.Range("C10").Value = 0.992162242584526
$ws.Range("E10").Value = 0.9967047873794244
$ws.Range("F10").Value = 0.9878110444644362
$ws.Range("J10").Value = 0.9988121622532352
$ws.Range("L10").Value = 1.000296919877218
$ws.Range("M10").Value = 0.9914381395819495
$ws.Range("N10").Value = 1.003274668619613
$ws.Range("C11").Value = 0.9906535521406713
$ws.Range("E11").Value = 0.9953936374089642
$ws.Range("F11").Value = 0.9858403664321077
$ws.Range("J11").Value = 0.9975575053047276
$ws.Range("L11").Value = 0.9991181050659222
$ws.Range("M11").Value = 0.989603955915083
$ws.Range("N11").Value = 1.002825974291096
$ws.Range("C12").Value = 0.9900913939142383
$ws.Range("E12").Value = 0.9949050722402972
$ws.Range("F12").Value = 0.9851061820500314
$ws.Range("J12").Value = 0.9970897735031711
$ws.Range("L12").Value = 0.9986786442991346
$ws.Range("M12").Value = 0.9889204816658861
$ws.Range("N12").Value = 1.002658519692851
$ws.Range("C13").Value = 0.9902120595381853
$ws.Range("E13").Value = 0.9950099418895324
$ws.Range("F13").Value = 0.9852637674690582
$ws.Range("J13").Value = 0.9971901811635614
$ws.Range("L13").Value = 0.9987729831400926
$ws.Range("M13").Value = 0.9890671890455454
$ws.Range("N13").Value = 1.00269447529006
$ws.Range("C14").Value = 0.9906071202911291
$ws.Range("E14").Value = 0.9953532842849069
$ws.Range("F14").Value = 0.9857797236426981
$ws.Range("J14").Value = 0.9975188773544151
$ws.Range("L14").Value = 0.9990818119466083
$ws.Range("M14").Value = 0.9895475046040828
$ws.Range("N14").Value = 1.002812148652134
$ws.Range("C15").Value = 0.9908502947745093
$ws.Range("E15").Value = 0.9955646225320273
$ws.Range("F15").Value = 0.9860973288618587
$ws.Range("J15").Value = 0.9977211715879972
$ws.Range("L15").Value = 0.9992718785883342
$ws.Range("M15").Value = 0.9898431519607187
$ws.Range("N15").Value = 1.002884545935585
$ws.Range("C16").Value = 0.9922621254230287
$ws.Range("E16").Value = 0.9967915900936594
$ws.Range("F16").Value = 0.9879415300188071
$ws.Range("J16").Value = 0.9988951946573905
$ws.Range("L16").Value = 1.00037493267123
$ws.Range("M16").Value = 0.9915595675619318
$ws.Range("N16").Value = 1.003304337241511
$ws.Range("C17").Value = 0.9931446518951078
$ws.Range("E17").Value = 0.9975585341355971
$ws.Range("F17").Value = 0.9890945458343663
$ws.Range("J17").Value = 0.9996286609749276
$ws.Range("L17").Value = 1.001064055618971
$ws.Range("M17").Value = 0.9926324388772686
$ws.Range("N17").Value = 1.003566273402792
$ws.Range("C18").Value = 0.9936583247212826
$ws.Range("E18").Value = 0.99800492249386
$ws.Range("F18").Value = 0.989765737786916
$ws.Range("J18").Value = 1.000055426526088
$ws.Range("L18").Value = 1.001465017693937
$ws.Range("M18").Value = 0.9932568864796906
$ws.Range("N18").Value = 1.003718561842482
$ws.Range("C19").Value = 0.9938332907999206
$ws.Range("E19").Value = 0.9981569685386651
$ws.Range("F19").Value = 0.9899943717553472
$ws.Range("J19").Value = 1.000200765487867
$ws.Range("L19").Value = 1.001601568606567
$ws.Range("M19").Value = 0.9934695821076406
$ws.Range("N19").Value = 1.003770404919007
$ws.Range("C20").Value = 0.9930500783226618
$ws.Range("E20").Value = 0.9974763476563394
$ws.Range("F20").Value = 0.9889709776493112
$ws.Range("J20").Value = 0.9995500762588494
$ws.Range("L20").Value = 1.000990222161191
$ws.Range("M20").Value = 0.9925174691326856
$ws.Range("N20").Value = 1.003538221405817
$ws.Range("C21").Value = 0.9904908339126369
$ws.Range("E21").Value = 0.995252221582819
$ws.Range("F21").Value = 0.9856278484838733
$ws.Range("J21").Value = 0.9974221317641369
$ws.Range("L21").Value = 0.9989909140153372
$ws.Range("M21").Value = 0.989406124396305
$ws.Range("N21").Value = 1.002777518719308
$ws.Range("C22").Value = 0.9888715021031139
$ws.Range("E22").Value = 0.9938448567504572
$ws.Range("F22").Value = 0.9835131937094163
$ws.Range("J22").Value = 0.9960743677744535
$ws.Range("L22").Value = 0.9977246088080325
$ws.Range("M22").Value = 0.9874372626241341
$ws.Range("N22").Value = 1.002294657611213
$ws.Range("C23").Value = 0.989730930884454
$ws.Range("E23").Value = 0.9945917945827535
$ws.Range("F23").Value = 0.9846354452697472
$ws.Range("J23").Value = 0.9967897932096137
$ws.Range("L23").Value = 0.998396795041804
$ws.Range("M23").Value = 0.9884822189254502
$ws.Range("N23").Value = 1.002551071293827
$ws.Range("C24").Value = 0.9930928154140612
$ws.Range("E24").Value = 0.9975134871437004
$ws.Range("F24").Value = 0.9890268169374213
$ws.Range("J24").Value = 0.9995855885604163
$ws.Range("L24").Value = 1.001023587384582
$ws.Range("M24").Value = 0.9925694231479878
$ws.Range("N24").Value = 1.003550898422906
$ws.Range("C25").Value = 0.9969515404821372
$ws.Range("E25").Value = 1.00086657049586
$ws.Range("F25").Value = 0.9940703804778079
$ws.Range("J25").Value = 1.002788820867279
$ws.Range("L25").Value = 1.004033090616584
$ws.Range("M25").Value = 0.9972601003021525
$ws.Range("N25").Value = 1.004691803431707
